$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "51.558.72"
$ws.Range("E2").Value = "  -0.62%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.942.19"
$ws.Range("E3").Value = "  +0.63%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "350.25"
$ws.Range("E5").Value = "  -0.59%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.33"
$ws.Range("E6").Value = "  -5.25%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -1.16%  "

# Row 9 - Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.603"
$ws.Range("E9").Value = "  -2.89%  "

# Row 10 - Avalanche
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.77"
$ws.Range("E10").Value = "  -4.01%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +1.27%  "

# Row 12 - Dogecoin
$ws.Range("E12").Value = "  -3.63%  "

# Row 13 - Chainlink
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.95"
$ws.Range("E13").Value = "  -5.72%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "3.394.27"
$ws.Range("E14").Value = "  +0.28%  "

# Row 15 - Polkadot
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.57"
$ws.Range("E15").Value = "  -2.49%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "2.933.73"
$ws.Range("E16").Value = "  +0.27%  "

# Row 17 - Polygon
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.965"
$ws.Range("E17").Value = "  -1.83%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "51.509.09"
$ws.Range("E18").Value = "  -0.79%  "

# Row 19 - ImmutableX
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.37"
$ws.Range("E19").Value = "  +1.76%  "

# Row 20 - Uniswap
$ws.Range("E20").Value = "  -3.02%  "

# Row 21 - InternetComputer(DFINITY)
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.40"
$ws.Range("E21").Value = "  -5.82%  "

# Row 22 - ShibaInu
$ws.Range("E22").Value = "  -1.61%  "

# Row 23 - Litecoin
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.96"
$ws.Range("E23").Value = "  -3.14%  "

# Row 24 - BitcoinCash
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "260.63"
$ws.Range("E24").Value = "  -2.91%  "

# Row 25 - PancakeSwap
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.70"
$ws.Range("E25").Value = "  -2.87%  "

# Row 26 - Kaspa
$ws.Range("E26").Value = "  -3.56%  "

# Row 27 - EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.45"
$ws.Range("E27").Value = "  -1.82%  "

# Row 28 - Dai
$ws.Range("E28").Value = "  +0.08%  "

# Row 29 - Filecoin
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.26"
$ws.Range("E29").Value = "  -0.09%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  +0.41%  "

# Row 31 - Cosmos
$ws.Range("E31").Value = "  -3.57%  "

# Row 32 & 33 - RenderToken and Toncoin swap positions
$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.19"
$ws.Range("E32").Value = "  -2.72%  "

$ws.Range("B33").Value = "RenderToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.05"
$ws.Range("E33").Value = "  -2.70%  "

# Row 34 - InjectiveProtocol
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "35.61"
$ws.Range("E34").Value = "  -4.28%  "

# Row 35 - OKB
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.47"
$ws.Range("E35").Value = "  -4.58%  "

# Row 36 - VeChain
$ws.Range("E36").Value = "  -5.24%  "

# Row 37 - FirstDigitalUSD
$ws.Range("E37").Value = "  -0.07%  "

# Row 38 - LidoDAOToken
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.13"
$ws.Range("E38").Value = "  -6.11%  "

# Row 39 - Celestia
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.65"
$ws.Range("E39").Value = "  -5.40%  "

# Row 40 - ARBITRUM
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.94"

# Row 41 - Stacks
$ws.Range("E41").Value = "  -1.52%  "

# Row 42 - Stellar
$ws.Range("E42").Value = "  -1.94%  "

# Row 43 - Monero
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "123.29"
$ws.Range("E43").Value = "  +10.94%  "

# Row 44 - EnergySwap
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.16"
$ws.Range("E44").Value = "  -5.28%  "

# Row 45 - WEMIXToken
$ws.Range("E45").Value = "  -3.33%  "

# Row 46 - Maker
$ws.Range("D46").Value = "2.099.35"
$ws.Range("E46").Value = "  -3.33%  "

# Row 47 - NEARProtocol
$ws.Range("E47").Value = "  -5.94%  "

# Row 48 - ApeXProtocol
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.30"
$ws.Range("E48").Value = "  -9.17%  "

# Row 49 - TheGraph
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.235"
$ws.Range("E49").Value = "  -5.32%  "

# Row 50 - BEAM
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0332"
$ws.Range("E50").Value = "  -4.48%  "

# Row 51 - SEI
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.906"
$ws.Range("E51").Value = "  -4.30%  "
